# "added player choosing reward screen"
#
# The localization sheet gets one new key/value pair appended as a new row
# (row 59): FORMATTED_UI_REWARDS_PLAYER_CHOOSING / "{0} is choosing rewards",
# with the other language columns (it/es/el) left as the "XXXX"
# not-yet-translated placeholder, matching every other still-untranslated
# row (e.g. rows 56-58). The view also scrolled down and the selection
# moved to B60 (just past the newly added row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 58 is the closest existing row with the exact look of the new row
# (key in col A, english text in col B, "XXXX" placeholders in C:E, and the
# unstyled/default row height). Copy it down to row 59 so the new row picks
# up the identical cell styles (s="3", s="13", s="4", s="4", s="13"), then
# overwrite the two real values.
$ws.Rows.Item(58).Copy()
$ws.Rows.Item(59).PasteSpecial()

$ws.Range("A59").Value2 = "FORMATTED_UI_REWARDS_PLAYER_CHOOSING"
$ws.Range("B59").Value2 = "{0} is choosing rewards"
$ws.Range("C59").Value2 = "XXXX"
$ws.Range("D59").Value2 = "XXXX"
$ws.Range("E59").Value2 = "XXXX"

# Clear the clipboard marching ants left over from the row copy.
$excel.CutCopyMode = 0

# Reflect the author scrolling down to/selecting just past the new row.
$ws.Range("B60").Select()
